# Update "想去人数" (interest count) and a couple of "最低票价" (min price)
# figures across the 展览 / 演出 / 全部类型 sheets to match the refreshed
# data snapshot generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibitions) ---
$ws1.Range("F2").Value = 29
$ws1.Range("F4").Value = 818
$ws1.Range("F5").Value = 868
$ws1.Range("F6").Value = 705
$ws1.Range("F7").Value = 1263
$ws1.Range("F9").Value = 870
$ws1.Range("F10").Value = 722
$ws1.Range("F13").Value = 384
$ws1.Range("F15").Value = 1030
$ws1.Range("F16").Value = 11390
$ws1.Range("F17").Value = 657
$ws1.Range("F22").Value = 291
$ws1.Range("F23").Value = 1803
$ws1.Range("F26").Value = 499
$ws1.Range("F28").Value = 110
$ws1.Range("F29").Value = 304
$ws1.Range("F31").Value = 269
$ws1.Range("F32").Value = 85
$ws1.Range("F33").Value = 105
$ws1.Range("F35").Value = 188
$ws1.Range("F36").Value = 207
$ws1.Range("F37").Value = 1204

# --- 演出 (Performances) ---
$ws2.Range("F7").Value = 153
$ws2.Range("G7").Value = 78
$ws2.Range("F12").Value = 97
$ws2.Range("F16").Value = 327

# --- 全部类型 (All types, combined view) ---
$ws4.Range("F3").Value = 29
$ws4.Range("F5").Value = 818
$ws4.Range("F7").Value = 868
$ws4.Range("F8").Value = 705
$ws4.Range("F9").Value = 1263
$ws4.Range("F12").Value = 153
$ws4.Range("G12").Value = 78
$ws4.Range("F13").Value = 870
$ws4.Range("F14").Value = 722
$ws4.Range("F17").Value = 1030
$ws4.Range("F18").Value = 11390
$ws4.Range("F20").Value = 657
$ws4.Range("F23").Value = 291
$ws4.Range("F24").Value = 1803
$ws4.Range("F26").Value = 499
$ws4.Range("F28").Value = 97
$ws4.Range("F29").Value = 97
$ws4.Range("F33").Value = 327
$ws4.Range("F34").Value = 304
$ws4.Range("F37").Value = 269
$ws4.Range("F38").Value = 85
$ws4.Range("F39").Value = 105
$ws4.Range("F42").Value = 188
$ws4.Range("F45").Value = 207
$ws4.Range("F46").Value = 1204
